$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# 1) The two existing worker rows (16: JAIR .. / 17: MIRTHA ..) swap places
#    so MIRTHA now comes first (row 16) and JAIR second (row 17).
# 2) A brand-new row is inserted below them (new row 18) for JAIR with an
#    extra "Periodo Mora" (2508), reusing his existing Valor Mora / Salario
#    Basico amounts.
# 3) The totals at the top of the statement are refreshed to reflect the new
#    row: VALOR MORA (E11) and Cant. Periodos (F13).
# 4) The trailing signature block (rows 22-23) shifts down one row (23-24)
#    to make room for the new data row.
# ---------------------------------------------------------------------------

# Insert a fresh row right after the current last data row (17). This pushes
# the blank spacer rows and the signature block (old rows 22-23) down to
# rows 23-24.
$ws.Rows("18").Insert()

# Row 17 currently still holds the original MIRTHA row (with its "closing"
# border style, i.e. the style used for the last row of the table). Move
# that whole row - formatting and all - down into the new row 18 first.
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))

# Row 16 holds the original JAIR row using the "middle" row border style.
# Duplicate that formatting+content into row 17, which now becomes a normal
# (non-closing) row of the table.
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))

# --- Row 16: now holds MIRTHA ENILSA MENDOZA SOTO ---------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "20220753"
$ws.Range("D16").Value = "MIRTHA ENILSA MENDOZA SOTO"
$ws.Range("E16").Value = "1908"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 828116

# --- Row 17: now holds JAIR JOSE FONTANILLA ARRIETA (period 2507) ----------
# (values already correct from the row-16 copy above - JAIR's original data -
# so nothing else to change here.)

# --- Row 18 (new, closing style from the old MIRTHA row): JAIR JOSE
#     FONTANILLA ARRIETA again, but for the new period 2508 ----------------
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "19874588"
$ws.Range("D18").Value = "JAIR JOSE FONTANILLA ARRIETA"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 72940
$ws.Range("G18").Value = 1823500

# --- Refresh the summary totals ---------------------------------------------
$ws.Range("E11").Value = 179005
$ws.Range("F13").Value = 3
